$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2370
$ws.Range("I62").Value = 1712.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 1712.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1088.5
$ws.Range("N62").Value = -6248
$ws.Range("H64").Value = 111116400
$ws.Range("I64").Value = 4457.4
$ws.Range("K64").Value = 4457.4
$ws.Range("M64").Value = -4209.4
$ws.Range("H65").Value = 2370
$ws.Range("I65").Value = 1712.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 8562.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -5442.5
$ws.Range("N65").Value = -31240
$ws.Range("H67").Value = 111116400
$ws.Range("I67").Value = 4457.4
$ws.Range("K67").Value = 4457.4
$ws.Range("M67").Value = -3599.4
$ws.Range("H86").Value = 8289.5
$ws.Range("I86").Value = 7696.2
$ws.Range("J86").Value = 8713.286
$ws.Range("K86").Value = 7696.2
$ws.Range("L86").Value = 8713.286
$ws.Range("M86").Value = -6573.2
$ws.Range("N86").Value = -10959.286
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 8289.5
$ws.Range("I89").Value = 7696.2
$ws.Range("J89").Value = 8713.286
$ws.Range("K89").Value = 38481
$ws.Range("L89").Value = 43566.43
$ws.Range("M89").Value = -32865
$ws.Range("N89").Value = -54798.43
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H98").Value = 4865.5557
$ws.Range("I98").Value = 5004.7646
$ws.Range("J98").Value = 2499
$ws.Range("K98").Value = 5004.7646
$ws.Range("L98").Value = 2499
$ws.Range("M98").Value = -3506.7646
$ws.Range("N98").Value = -5495
$ws.Range("H103").Value = 2336.5
$ws.Range("I103").Value = 2541
$ws.Range("J103").Value = 1927.5
$ws.Range("K103").Value = 7623
$ws.Range("L103").Value = 5782.5
$ws.Range("M103").Value = -7037
$ws.Range("N103").Value = -6954.5
$ws.Range("H122").Value = 4865.5557
$ws.Range("I122").Value = 5004.7646
$ws.Range("J122").Value = 2499
$ws.Range("K122").Value = 15014.2938
$ws.Range("L122").Value = 7497
$ws.Range("M122").Value = -12564.2938
$ws.Range("N122").Value = -12397
$ws.Range("H135").Value = 1044.2858
$ws.Range("I135").Value = 916.8182
$ws.Range("J135").Value = 1511.6666
$ws.Range("K135").Value = 8251.363800000001
$ws.Range("L135").Value = 13604.9994
$ws.Range("M135").Value = -5716.363800000001
$ws.Range("N135").Value = -18674.9994
$ws.Range("H138").Value = 3396.4722
$ws.Range("J138").Value = 2989.8462
$ws.Range("L138").Value = 8969.5386
$ws.Range("N138").Value = -19249.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1615.1632
$ws.Range("I32").Value = 1637.3405
$ws.Range("K32").Value = 1637.3405
$ws.Range("M32").Value = -1350.3405
$ws.Range("H61").Value = 2759.1177
$ws.Range("I61").Value = 2051.25
$ws.Range("K61").Value = 2051.25
$ws.Range("M61").Value = -1839.25
$ws.Range("H102").Value = 3709.0667
$ws.Range("I102").Value = 2828.1667
$ws.Range("J102").Value = 7232.6665
$ws.Range("K102").Value = 2828.1667
$ws.Range("L102").Value = 7232.6665
$ws.Range("M102").Value = -1206.1667
$ws.Range("N102").Value = -10476.6665
$ws.Range("H136").Value = 2759.1177
$ws.Range("I136").Value = 2051.25
$ws.Range("K136").Value = 6153.75
$ws.Range("M136").Value = -3603.75
$ws.Range("H141").Value = 54199.4
$ws.Range("J141").Value = 54199.4
$ws.Range("L141").Value = 54199.4
$ws.Range("N141").Value = -64559.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1003.625
$ws.Range("I64").Value = 720.75
$ws.Range("K64").Value = 720.75
$ws.Range("M64").Value = -495.75
$ws.Range("H67").Value = 1003.625
$ws.Range("I67").Value = 720.75
$ws.Range("K67").Value = 720.75
$ws.Range("M67").Value = 59.25
$ws.Range("H99").Value = 2778.8572
$ws.Range("J99").Value = 5122.75
$ws.Range("L99").Value = 5122.75
$ws.Range("N99").Value = -8118.75
$ws.Range("H105").Value = 65002000
$ws.Range("I105").Value = 5001999.5
$ws.Range("J105").Value = 125002000
$ws.Range("K105").Value = 5001999.5
$ws.Range("L105").Value = 125002000
$ws.Range("M105").Value = -5000252.5
$ws.Range("N105").Value = -125005494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4094.5532
$ws.Range("I31").Value = 4615.0713
$ws.Range("J31").Value = 3873.7273
$ws.Range("K31").Value = 4615.0713
$ws.Range("L31").Value = 3873.7273
$ws.Range("M31").Value = -4320.0713
$ws.Range("N31").Value = -4463.7273
$ws.Range("H34").Value = 4094.5532
$ws.Range("I34").Value = 4615.0713
$ws.Range("J34").Value = 3873.7273
$ws.Range("K34").Value = 4615.0713
$ws.Range("L34").Value = 3873.7273
$ws.Range("M34").Value = -4413.0713
$ws.Range("N34").Value = -4277.7273
$ws.Range("H102").Value = 28399.5
$ws.Range("J102").Value = 28399.5
$ws.Range("L102").Value = 28399.5
$ws.Range("N102").Value = -33267.5
$ws.Range("H132").Value = 5559640.5
$ws.Range("I132").Value = 3928.7292
$ws.Range("K132").Value = 11786.1876
$ws.Range("M132").Value = -9256.187600000001
$ws.Range("H134").Value = 2693.2559
$ws.Range("J134").Value = 4389
$ws.Range("L134").Value = 13167
$ws.Range("N134").Value = -18237

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I82").Value = 11500
$ws.Range("K82").Value = 34500
$ws.Range("M82").Value = -34094
$ws.Range("I85").Value = 11500
$ws.Range("K85").Value = 34500
$ws.Range("M85").Value = -33096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3597.2727
$ws.Range("I122").Value = 2902.2
$ws.Range("K122").Value = 8706.599999999999
$ws.Range("M122").Value = -6256.599999999999
$ws.Range("H141").Value = 69485.39999999999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 69485.39999999999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 69485.39999999999
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -79845.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H122").Value = 12832.833
$ws.Range("I122").Value = 8600
$ws.Range("K122").Value = 25800
$ws.Range("M122").Value = -23350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I100").Value = 1000000000
$ws.Range("K100").Value = 2000000000
$ws.Range("M100").Value = -1999999459
$ws.Range("H132").Value = 3916.5862
$ws.Range("I132").Value = 3817.08
$ws.Range("J132").Value = 4538.5
$ws.Range("K132").Value = 11451.24
$ws.Range("L132").Value = 13615.5
$ws.Range("M132").Value = -8921.24
$ws.Range("N132").Value = -18675.5
